# Insert a new daily price record at row 563 of the "Poroto verde" sheet.
# All existing rows from 563 downward shift down by one (563->564, ..., 635->636),
# and the freshly inserted row 563 is populated with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 563..635 down by one row, leaving a blank row 563 behind.
$ws.Rows.Item(563).Insert()

# Populate the newly inserted row 563 with the new data point.
$ws.Range("A563").Value = 3
$ws.Range("B563").Value = "Femacal de La Calera"
$ws.Range("C563").Value = "Coquimbo"
$ws.Range("D563").Value = 45124
$ws.Range("E563").Value = 5
$ws.Range("F563").Value = 100112031
$ws.Range("G563").Value = "Poroto verde"
$ws.Range("H563").Value = "Magnum"
$ws.Range("I563").Value = "Primera"
$ws.Range("J563").Value = 70
$ws.Range("K563").Value = 25000
$ws.Range("L563").Value = 26000
$ws.Range("M563").Value = 25500
$ws.Range("N563").Value = "`$/malla 25 kilos"
$ws.Range("O563").Value = "Región de Arica y Parinacota"
$ws.Range("P563").Value = 1020
$ws.Range("Q563").Value = 25
$ws.Range("R563").Value = "Hortaliza"
